$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values are plain numeric-looking text (e.g. "572.36").
# A direct .Value assignment would let Excel auto-convert these to real numbers,
# which both loses the exact text ("0.860" -> 0.86) and changes the cell type.
# Forcing a Text number format before the assignment keeps them as literal text;
# re-applying the Normal style afterwards removes the only formatting residue so
# the cell style matches the untouched cells around it.
function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "64.284.59"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.428.48"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "572.36"
$ws.Range("E5").Value = "  -0.31%  "
Set-TextValue "D6" "161.80"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.428.68"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -8.47%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  -2.33%  "
Set-TextValue "D12" "0.426"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "4.017.40"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  -1.73%  "
Set-TextValue "D16" "0.0000175"
$ws.Range("E16").Value = "  -7.11%  "
$ws.Range("D17").Value = "64.340.69"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "3.410.73"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("E19").Value = "  -3.86%  "
Set-TextValue "D20" "13.65"
$ws.Range("E20").Value = "  -1.70%  "
Set-TextValue "D21" "379.86"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("E22").Value = "  -1.61%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.47%  "
Set-TextValue "D24" "71.65"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -5.50%  "
$ws.Range("E26").Value = "  -1.91%  "
Set-TextValue "D27" "9.67"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  +0.09%  "
Set-TextValue "D30" "6.09"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("E32").Value = "  -0.21%  "
Set-TextValue "D33" "23.06"
$ws.Range("E33").Value = "  -1.03%  "
Set-TextValue "D34" "7.10"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -5.28%  "
Set-TextValue "D36" "159.08"
$ws.Range("E36").Value = "  -0.78%  "
Set-TextValue "D37" "0.860"
$ws.Range("E37").Value = "  +11.22%  "
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.818.60"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.0733"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  -2.67%  "
Set-TextValue "D42" "43.22"
Set-TextValue "D43" "25.88"
$ws.Range("E43").Value = "  -2.30%  "
Set-TextValue "D44" "26.56"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -3.42%  "
Set-TextValue "D47" "346.79"
$ws.Range("E47").Value = "  +8.99%  "
Set-TextValue "D48" "2.41"
$ws.Range("E48").Value = "  +6.86%  "
Set-TextValue "D49" "1.06"
$ws.Range("E49").Value = "  -1.24%  "
Set-TextValue "D50" "6.34"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("E51").Value = "  -4.27%  "
